$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and G to be treated as text so numeric-looking strings
# (prices, and the "0" hour value) are stored as text, matching the source format.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '246.36'
$ws.Range("F2").Value = '30-12-2022'
$ws.Range("G2").Value = '0'
$ws.Range("D3").Value = '24.29'
$ws.Range("F3").Value = '30-12-2022'
$ws.Range("G3").Value = '0'
$ws.Range("D4").Value = '5.278'
$ws.Range("F4").Value = '30-12-2022'
$ws.Range("G4").Value = '0'
$ws.Range("D5").Value = '0.05812'
$ws.Range("F5").Value = '30-12-2022'
$ws.Range("G5").Value = '0'
$ws.Range("D6").Value = '6.517'
$ws.Range("F6").Value = '30-12-2022'
$ws.Range("G6").Value = '0'
$ws.Range("D7").Value = '3.143'
$ws.Range("F7").Value = '30-12-2022'
$ws.Range("G7").Value = '0'
$ws.Range("D8").Value = '0.8167'
$ws.Range("F8").Value = '30-12-2022'
$ws.Range("G8").Value = '0'
$ws.Range("D9").Value = '0.8536'
$ws.Range("F9").Value = '30-12-2022'
$ws.Range("G9").Value = '0'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '0.0005957'
$ws.Range("E10").Value = '9OneONE'
$ws.Range("F10").Value = '30-12-2022'
$ws.Range("G10").Value = '0'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1363'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("F11").Value = '30-12-2022'
$ws.Range("G11").Value = '0'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.06941'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("F12").Value = '30-12-2022'
$ws.Range("G12").Value = '0'
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = '0.03127'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("F13").Value = '30-12-2022'
$ws.Range("G13").Value = '0'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '0.02875'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("F14").Value = '30-12-2022'
$ws.Range("G14").Value = '0'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '0.09402'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("F15").Value = '30-12-2022'
$ws.Range("G15").Value = '0'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = '3.753'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("F16").Value = '30-12-2022'
$ws.Range("G16").Value = '0'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = '0.001513'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("F17").Value = '30-12-2022'
$ws.Range("G17").Value = '0'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = '0.04676'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("F18").Value = '30-12-2022'
$ws.Range("G18").Value = '0'
$ws.Range("D19").Value = '0.006264'
$ws.Range("F19").Value = '30-12-2022'
$ws.Range("G19").Value = '0'
$ws.Range("D20").Value = '0.001234'
$ws.Range("F20").Value = '30-12-2022'
$ws.Range("G20").Value = '0'
$ws.Range("D21").Value = '0.004623'
$ws.Range("F21").Value = '30-12-2022'
$ws.Range("G21").Value = '0'
$ws.Range("D22").Value = '0.00006897'
$ws.Range("E22").Value = '21NitroExNTXWorstin24h'
$ws.Range("F22").Value = '30-12-2022'
$ws.Range("G22").Value = '0'
$ws.Range("D23").Value = '3.499'
$ws.Range("F23").Value = '30-12-2022'
$ws.Range("G23").Value = '0'
$ws.Range("D24").Value = '2.148'
$ws.Range("F24").Value = '30-12-2022'
$ws.Range("G24").Value = '0'
$ws.Range("D25").Value = '0.3193'
$ws.Range("F25").Value = '30-12-2022'
$ws.Range("G25").Value = '0'
$ws.Range("D26").Value = '0.1337'
$ws.Range("F26").Value = '30-12-2022'
$ws.Range("G26").Value = '0'
$ws.Range("F27").Value = '30-12-2022'
$ws.Range("G27").Value = '0'
$ws.Range("D28").Value = '0.0002328'
$ws.Range("F28").Value = '30-12-2022'
$ws.Range("G28").Value = '0'
$ws.Range("F29").Value = '30-12-2022'
$ws.Range("G29").Value = '0'
$ws.Range("F30").Value = '30-12-2022'
$ws.Range("G30").Value = '0'
$ws.Range("F31").Value = '30-12-2022'
$ws.Range("G31").Value = '0'
$ws.Range("F32").Value = '30-12-2022'
$ws.Range("G32").Value = '0'
$ws.Range("F33").Value = '30-12-2022'
$ws.Range("G33").Value = '0'
$ws.Range("F34").Value = '30-12-2022'
$ws.Range("G34").Value = '0'
$ws.Range("F35").Value = '30-12-2022'
$ws.Range("G35").Value = '0'
$ws.Range("F36").Value = '30-12-2022'
$ws.Range("G36").Value = '0'
$ws.Range("F37").Value = '30-12-2022'
$ws.Range("G37").Value = '0'
$ws.Range("F38").Value = '30-12-2022'
$ws.Range("G38").Value = '0'
$ws.Range("F39").Value = '30-12-2022'
$ws.Range("G39").Value = '0'
$ws.Range("D40").Value = '0.03672'
$ws.Range("F40").Value = '30-12-2022'
$ws.Range("G40").Value = '0'
$ws.Range("D41").Value = '0.006237'
$ws.Range("F41").Value = '30-12-2022'
$ws.Range("G41").Value = '0'
$ws.Range("F42").Value = '30-12-2022'
$ws.Range("G42").Value = '0'
$ws.Range("F43").Value = '30-12-2022'
$ws.Range("G43").Value = '0'
$ws.Range("D44").Value = '0.007488'
$ws.Range("F44").Value = '30-12-2022'
$ws.Range("G44").Value = '0'
$ws.Range("D45").Value = '0.00005255'
$ws.Range("F45").Value = '30-12-2022'
$ws.Range("G45").Value = '0'
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("F46").Value = '30-12-2022'
$ws.Range("G46").Value = '0'
$ws.Range("D47").Value = '0.3698'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("F47").Value = '30-12-2022'
$ws.Range("G47").Value = '0'
$ws.Range("D48").Value = '0.002403'
$ws.Range("F48").Value = '30-12-2022'
$ws.Range("G48").Value = '0'
$ws.Range("F49").Value = '30-12-2022'
$ws.Range("G49").Value = '0'
$ws.Range("F50").Value = '30-12-2022'
$ws.Range("G50").Value = '0'
$ws.Range("F51").Value = '30-12-2022'
$ws.Range("G51").Value = '0'
